# The underlying data rows for this species-observation sheet (rows 3-7)
# got reshuffled: the per-observation fields (columns A, B, D, E, F, G, H,
# M, Q, R) were redistributed among rows 3-7, while all the shared
# metadata columns (C, I-L, N, P, S-Z, AA, AB, AD-AG, AT, AW-AY) stayed
# untouched (they already held identical values across these rows).
#
# Permutation of the moving columns, expressed as target-row <- source-row
# (captured from the before-workbook state, since we must not read a row
# after it has already been overwritten):
#   row 3 <- row 4
#   row 4 <- row 6
#   row 5 <- row 7
#   row 6 <- row 3
#   row 7 <- row 5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$movingCols = @("A", "B", "D", "E", "F", "G", "H", "M", "Q", "R")
$mapping = @{ 3 = 4; 4 = 6; 5 = 7; 6 = 3; 7 = 5 }

# Snapshot the "moving" cell values for every source row before any writes
# happen, so later writes never clobber data we still need to read.
$snapshot = @{}
foreach ($row in 3..7) {
    $rowData = @{}
    foreach ($col in $movingCols) {
        $rowData[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowData
}

foreach ($targetRow in $mapping.Keys) {
    $sourceRow = $mapping[$targetRow]
    $rowData = $snapshot[$sourceRow]
    foreach ($col in $movingCols) {
        $val = $rowData[$col]
        if ($null -eq $val -or $val -eq "") {
            $ws.Range("$col$targetRow").ClearContents()
        } else {
            $ws.Range("$col$targetRow").Value = $val
        }
    }
}
